# Scheduled runner update: refresh market price / profit figures across all Sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 152.625
$ws.Range("I61").Value = 152.625
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 457.875
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -285.875
$ws.Range("N61").ClearContents()
$ws.Range("H64").Value = 4188.3706
$ws.Range("I64").Value = 4084.2942
$ws.Range("J64").Value = 4365.3
$ws.Range("K64").Value = 4084.2942
$ws.Range("L64").Value = 4365.3
$ws.Range("M64").Value = -3836.2942
$ws.Range("N64").Value = -4861.3
$ws.Range("H67").Value = 4188.3706
$ws.Range("I67").Value = 4084.2942
$ws.Range("J67").Value = 4365.3
$ws.Range("K67").Value = 4084.2942
$ws.Range("L67").Value = 4365.3
$ws.Range("M67").Value = -3226.2942
$ws.Range("N67").Value = -6081.3
$ws.Range("H92").Value = 1311.3
$ws.Range("J92").Value = 1105.6666
$ws.Range("L92").Value = 1105.6666
$ws.Range("N92").Value = -3601.6666
$ws.Range("H113").Value = 2013.3334
$ws.Range("I113").Value = 2493.3333
$ws.Range("J113").Value = 1853.3334
$ws.Range("K113").Value = 2493.3333
$ws.Range("L113").Value = 1853.3334
$ws.Range("M113").Value = 760.6667000000002
$ws.Range("N113").Value = -8361.3334
$ws.Range("H129").Value = 826.0769
$ws.Range("J129").Value = 904.0303
$ws.Range("L129").Value = 2712.0909
$ws.Range("N129").Value = -12712.0909
$ws.Range("H137").Value = 40002412
$ws.Range("I137").Value = 1545
$ws.Range("K137").Value = 4635
$ws.Range("M137").Value = -2085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H32").Value = 4993.73
$ws.Range("I32").Value = 5003.7676
$ws.Range("K32").Value = 5003.7676
$ws.Range("M32").Value = -4716.7676
$ws.Range("H63").Value = 7999.5
$ws.Range("I63").Value = 7999
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 7999
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -7313
$ws.Range("N63").Value = -9372
$ws.Range("H66").Value = 7999.5
$ws.Range("I66").Value = 7999
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 39995
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -36563
$ws.Range("N66").Value = -46864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 247.93333
$ws.Range("I22").Value = 199.92857
$ws.Range("J22").Value = 920
$ws.Range("K22").Value = 199.92857
$ws.Range("L22").Value = 920
$ws.Range("M22").Value = -26.92857000000001
$ws.Range("N22").Value = -1266
$ws.Range("H105").Value = 1695.8823
$ws.Range("I105").Value = 1690
$ws.Range("J105").Value = 1740
$ws.Range("K105").Value = 1690
$ws.Range("L105").Value = 1740
$ws.Range("M105").Value = 57
$ws.Range("N105").Value = -5234

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1696.25
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 1928.3334
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1928.3334
$ws.Range("M8").Value = -860
$ws.Range("N8").Value = -2208.3334
$ws.Range("H62").Value = 2354.3333
$ws.Range("I62").Value = 2063.3333
$ws.Range("J62").Value = 3033.3333
$ws.Range("K62").Value = 2063.3333
$ws.Range("L62").Value = 3033.3333
$ws.Range("M62").Value = -1439.3333
$ws.Range("N62").Value = -4281.3333
$ws.Range("H65").Value = 2354.3333
$ws.Range("I65").Value = 2063.3333
$ws.Range("J65").Value = 3033.3333
$ws.Range("K65").Value = 10316.6665
$ws.Range("L65").Value = 15166.6665
$ws.Range("M65").Value = -7196.666499999999
$ws.Range("N65").Value = -21406.6665
$ws.Range("H68").Value = 18237
$ws.Range("J68").Value = 18237
$ws.Range("L68").Value = 18237
$ws.Range("N68").Value = -19735
$ws.Range("H71").Value = 18237
$ws.Range("J71").Value = 18237
$ws.Range("L71").Value = 54711
$ws.Range("N71").Value = -62199
$ws.Range("H134").Value = 3308.3635
$ws.Range("I134").Value = 3062.8125
$ws.Range("J134").Value = 3963.1667
$ws.Range("K134").Value = 9188.4375
$ws.Range("L134").Value = 11889.5001
$ws.Range("M134").Value = -6653.4375
$ws.Range("N134").Value = -16959.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 641.3333
$ws.Range("I113").Value = 666.6667
$ws.Range("J113").Value = 616
$ws.Range("K113").Value = 2000.0001
$ws.Range("L113").Value = 1848
$ws.Range("M113").Value = 169.9999
$ws.Range("N113").Value = -6188
$ws.Range("H122").Value = 13334201
$ws.Range("I122").Value = 16667161
$ws.Range("J122").Value = 2359
$ws.Range("K122").Value = 150004449
$ws.Range("L122").Value = 21231
$ws.Range("M122").Value = -150001999
$ws.Range("N122").Value = -26131
$ws.Range("H131").Value = 2754.75
$ws.Range("I131").Value = 17076.666
$ws.Range("J131").Value = 1877.898
$ws.Range("K131").Value = 51229.99800000001
$ws.Range("L131").Value = 5633.694
$ws.Range("M131").Value = -46189.99800000001
$ws.Range("N131").Value = -15713.694
$ws.Range("H137").Value = 22727.328
$ws.Range("J137").Value = 27028.512
$ws.Range("L137").Value = 81085.53599999999
$ws.Range("N137").Value = -91285.53599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7762.4194
$ws.Range("I70").Value = 10950.333
$ws.Range("J70").Value = 4773.75
$ws.Range("K70").Value = 10950.333
$ws.Range("L70").Value = 4773.75
$ws.Range("M70").Value = -10680.333
$ws.Range("N70").Value = -5313.75
$ws.Range("H73").Value = 7762.4194
$ws.Range("I73").Value = 10950.333
$ws.Range("J73").Value = 4773.75
$ws.Range("K73").Value = 10950.333
$ws.Range("L73").Value = 4773.75
$ws.Range("M73").Value = -10014.333
$ws.Range("N73").Value = -6645.75
$ws.Range("H80").Value = 3790
$ws.Range("I80").Value = 3487.5
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 3487.5
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -2489.5
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 3790
$ws.Range("I83").Value = 3487.5
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 17437.5
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -12445.5
$ws.Range("N83").Value = -34984
$ws.Range("H132").Value = 1662.619
$ws.Range("I132").Value = 1328.7142
$ws.Range("J132").Value = 1829.5714
$ws.Range("K132").Value = 3986.1426
$ws.Range("L132").Value = 5488.7142
$ws.Range("M132").Value = -1456.1426
$ws.Range("N132").Value = -10548.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 5015.2383
$ws.Range("J5").Value = 5015.2383
$ws.Range("L5").Value = 5015.2383
$ws.Range("N5").Value = -5241.2383
$ws.Range("H20").Value = 3285.7144
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5452
$ws.Range("H40").Value = 1616.3889
$ws.Range("I40").Value = 1610.3334
$ws.Range("K40").Value = 1610.3334
$ws.Range("M40").Value = -1474.3334
$ws.Range("H70").Value = 21500
$ws.Range("I70").Value = 14000
$ws.Range("J70").Value = 29000
$ws.Range("K70").Value = 14000
$ws.Range("L70").Value = 29000
$ws.Range("M70").Value = -13730
$ws.Range("N70").Value = -29540
$ws.Range("H73").Value = 21500
$ws.Range("I73").Value = 14000
$ws.Range("J73").Value = 29000
$ws.Range("K73").Value = 14000
$ws.Range("L73").Value = 29000
$ws.Range("M73").Value = -13064
$ws.Range("N73").Value = -30872
$ws.Range("H132").Value = 33626
$ws.Range("I132").Value = 43668
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 131004
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -128474
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3300003.8
$ws.Range("I14").Value = 3000000
$ws.Range("K14").Value = 3000000
$ws.Range("M14").Value = -2999832
$ws.Range("H20").Value = 87508.25
$ws.Range("I20").Value = 50000
$ws.Range("K20").Value = 50000
$ws.Range("M20").Value = -49760
$ws.Range("H24").Value = 100010
$ws.Range("J24").Value = 100010
$ws.Range("L24").Value = 100010
$ws.Range("N24").Value = -100470
$ws.Range("H30").Value = 10000
$ws.Range("J30").Value = 10000
$ws.Range("L30").Value = 10000
$ws.Range("N30").Value = -10214
$ws.Range("H122").Value = 1805.7778
$ws.Range("I122").Value = 1405.8182
$ws.Range("K122").Value = 4217.4546
$ws.Range("M122").Value = -1767.4546
$ws.Range("H132").Value = 3889.9814
$ws.Range("I132").Value = 5229.1113
$ws.Range("J132").Value = 1211.7222
$ws.Range("K132").Value = 15687.3339
$ws.Range("L132").Value = 3635.1666
$ws.Range("M132").Value = -13157.3339
$ws.Range("N132").Value = -8695.1666

Write-Output "Applied market data refresh across all sheets"
